$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column O to column P for rows 3-10 first
$ws.Range("O3:O10").Copy()
$ws.Range("P3:P10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("P4").Value = 2022
$ws.Range("P6").Value = 1373
$ws.Range("P7").Value = "-"
$ws.Range("P7").HorizontalAlignment = -4152
$ws.Range("P8").Value = 117
$ws.Range("P9").Value = 154
$ws.Range("P10").Value = 885

[void]$ws.Range("P7").Select()
